$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "fixed a bug with the 20 minute trade. Now you don't need to enter the
#  close price if the trader can't locate data from yahoo" -> a new trade
# row (row 4) is appended to the sheet, where the sell price (D) is an
# estimate/fallback rather than a real close price off Yahoo.

$ws.Range("A4").Value = 10015.91
$ws.Range("B4").Value = 10039
$ws.Range("C4").Value = 286.39
$ws.Range("D4").Value = 287.04000000000002
$ws.Range("E4").Value = $true
$ws.Range("F4").Value = 0.23
$ws.Range("H4").Value = $false

# G4 holds a date/time (same numeric date format as G3) -- copy the
# existing date cell's formatting over before writing the new value so the
# new cell reuses the workbook's existing date style rather than minting a
# new one.
$ws.Range("G3").Copy()
$ws.Range("G4").PasteSpecial(-4122)
$ws.Range("G4").Value = 42608.639652777776

# Column A's best-fit width shifts slightly (8.85546875 -> 9) now that the
# column holds the new row's data.
$ws.Columns("A").ColumnWidth = 8.166666666666666
